$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E3").Value = 16.621
$ws.Range("A4").Value = -20.824
$ws.Range("A6").Value = -22.082
$ws.Range("A7").Value = -21.255
$ws.Range("B7").Value = 6.536
$ws.Range("A8").Value = -21.561
$ws.Range("B11").Value = 6.169
$ws.Range("B12").Value = 5.481
$ws.Range("D12").Value = -6.825
$ws.Range("E12").Value = 17.344
$ws.Range("D13").Value = -8.115
$ws.Range("E13").Value = 16.593
$ws.Range("D14").Value = -7.661
$ws.Range("B15").Value = 5.359
$ws.Range("A16").Value = -21.485
$ws.Range("D16").Value = -8.610000000000001
$ws.Range("D19").Value = -7.796000000000001
$ws.Range("A20").Value = -22.103
$ws.Range("B20").Value = 5.928
$ws.Range("D20").Value = -7.856999999999999
$ws.Range("A21").Value = -20.952
$ws.Range("B21").Value = 7.692000000000002
$ws.Range("B22").Value = 6.557
$ws.Range("D22").Value = -8.151
$ws.Range("E22").Value = 16.492
$ws.Range("B23").Value = 7.398999999999999
$ws.Range("E25").Value = 17.331
$ws.Range("A28").Value = -21.778
$ws.Range("A29").Value = -21.312
$ws.Range("B29").Value = 6.093
$ws.Range("E29").Value = 17.143
$ws.Range("A30").Value = -21.744
$ws.Range("A32").Value = -21.663
$ws.Range("B34").Value = 7.603999999999999
$ws.Range("E34").Value = 16.665
$ws.Range("D36").Value = -7.632
$ws.Range("A40").Value = -20.407
$ws.Range("B42").Value = 7.047000000000001
$ws.Range("B43").Value = 5.635000000000001
$ws.Range("D43").Value = -7.99
$ws.Range("E43").Value = 16.9
$ws.Range("B44").Value = 5.339
$ws.Range("B45").Value = 5.403999999999999
$ws.Range("A46").Value = -20.849
$ws.Range("B46").Value = 6.865
$ws.Range("D46").Value = -8.059000000000001
$ws.Range("E48").Value = 17.347
$ws.Range("B50").Value = 5.637
$ws.Range("D50").Value = -7.897
$ws.Range("A51").Value = -20.952
$ws.Range("B51").Value = 7.930000000000001
$ws.Range("A52").Value = -21.588
$ws.Range("A57").Value = -21.64
$ws.Range("B57").Value = 6.046000000000001
$ws.Range("A59").Value = -22.145
$ws.Range("E60").Value = 16.369
$ws.Range("A62").Value = -21.645
$ws.Range("B65").Value = 5.245
$ws.Range("A66").Value = -21.345
$ws.Range("B66").Value = 5.895
$ws.Range("B67").Value = 5.897
$ws.Range("E68").Value = 17.409
$ws.Range("E70").Value = 17.642
$ws.Range("E71").Value = 17.409
$ws.Range("A73").Value = -20.618
$ws.Range("E73").Value = 16.507
$ws.Range("A74").Value = -21.043
$ws.Range("D76").Value = -8.015000000000001
$ws.Range("A77").Value = -21.298
$ws.Range("E78").Value = 16.559
$ws.Range("B79").Value = 5.605
$ws.Range("B84").Value = 5.799
$ws.Range("B87").Value = 4.447000000000001
$ws.Range("E87").Value = 16.585
$ws.Range("A92").Value = -21.103
$ws.Range("B92").Value = 6.044
$ws.Range("E92").Value = 17.789
$ws.Range("D95").Value = -7.921000000000001
$ws.Range("B97").Value = 5.161
$ws.Range("D97").Value = -8.465
$ws.Range("D99").Value = -7.772
$ws.Range("A100").Value = -21.481
$ws.Range("E101").Value = 16.556
